$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Find-ParagraphByText {
    param($doc, [string]$text, [bool]$exact = $true)
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $candidate = $doc.Paragraphs.Item($i)
        $t = $candidate.Range.Text.TrimEnd([char]13, [char]7)
        if ($exact) {
            if ($t -eq $text) { return $candidate }
        } else {
            if ($t.StartsWith($text)) { return $candidate }
        }
    }
    throw "Paragraph starting with '$text' not found"
}

# 1. Update the H1 title text ("Slot" removed before "Free").
$titlePara = $d.Paragraphs.Item(1)
$titleXml = '<w:p ' + $wNs + '><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Play Yu Tu Jin Cai Cash Collect Free | Review</w:t></w:r></w:p>'
$null = $titlePara.Range.InsertXML($titleXml)

# 2. Remove the whole "Meta description: ..." paragraph entirely (now paragraph 2).
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# 3. Update the "What we like" bullet list items (text only; formatting/pPr unchanged).
$listPPr = '<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr>'

$p = Find-ParagraphByText $d "Beautifully executed game symbols"
$xml = '<w:p ' + $wNs + '>' + $listPPr + '<w:r/><w:r><w:t>Simple and easy-to-understand gameplay mechanics</w:t></w:r></w:p>'
$null = $p.Range.InsertXML($xml)

$p = Find-ParagraphByText $d "Special Cash Collect feature for multipliers"
$xml = '<w:p ' + $wNs + '>' + $listPPr + '<w:r/><w:r><w:t>Stunning graphics that capture the essence of the slot</w:t></w:r></w:p>'
$null = $p.Range.InsertXML($xml)

$p = Find-ParagraphByText $d "Impressive four-tiered jackpot available"
$xml = '<w:p ' + $wNs + '>' + $listPPr + '<w:r/><w:r><w:t>Special features with multiplier and free spins</w:t></w:r></w:p>'
$null = $p.Range.InsertXML($xml)

$p = Find-ParagraphByText $d "Stunning Asian-inspired graphics"
$xml = '<w:p ' + $wNs + '>' + $listPPr + '<w:r/><w:r><w:t>Impressive jackpot tiers for big wins</w:t></w:r></w:p>'
$null = $p.Range.InsertXML($xml)

# 4. Update the "What we don't like" bullet list item.
$p = Find-ParagraphByText $d "Low RTP at only 94.9%"
$xml = '<w:p ' + $wNs + '>' + $listPPr + '<w:r/><w:r><w:t>RTP of 94.9% detracts from the overall experience</w:t></w:r></w:p>'
$null = $p.Range.InsertXML($xml)

# 5. Turn the "Not available for players in certain countries" bullet into a plain,
#    bold, style-less paragraph reading "Play Yu Tu Jin Cai Cash Collect Free | Review".
$p = Find-ParagraphByText $d "Not available for players in certain countries"
$xml = '<w:p ' + $wNs + '><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Yu Tu Jin Cai Cash Collect Free | Review</w:t></w:r></w:p>'
$null = $p.Range.InsertXML($xml)

# 6. Replace the "Prompt: ..." paragraph text with the new copy (keep the italics).
$p = Find-ParagraphByText $d "Prompt:" $false
$xml = '<w:p ' + $wNs + '><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Play Yu Tu Jin Cai Cash Collect for free and discover its stunning graphics and special features.</w:t></w:r></w:p>'
$null = $p.Range.InsertXML($xml)
